$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on Hoja1 (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$old = $wsHoja1.Range("A1").Value()
$new = $old.Replace(
    "1000 Bs = 7.35 = 30242.65 pesos",
    "1000 Bs = 7.33 = 30043.96 pesos"
)
$new = $new.Replace(
    "30242.65 pesos = 7.34 = 966.51 Bs",
    "30043.96 pesos = 7.29 = 961.62 Bs"
)
$wsHoja1.Range("A1").Value = $new

# --- Update the rate cells on the "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 136.5
$wsTasas.Range("O10").Value = 4101
$wsTasas.Range("N12").Value = 4120.95
$wsTasas.Range("O12").Value = 131.9
